$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 6: "Retângulo de cantos arredondados 10" (roundRect) ---
$roundRect = $s.Shapes.Item(6)
$roundRect.Left = 42.484409448818894
$roundRect.Top = 71.46622047244094

$rrTextRange = $roundRect.TextFrame.TextRange
$newRun = $rrTextRange.InsertAfter("auklaaau")
$newRun.LanguageID = "pt-BR"
$newRun.Font.Bold = $true

# --- Shape 7: "CaixaDeTexto 16" (DISCIPLINA textbox) ---
$discShape = $s.Shapes.Item(7)
$discShape.Left = 59.31094488188976
$discShape.Top = 83.80692913385826
$discShape.Width = 436.7671653543307
$discShape.Height = 50.892204724409446
$discShape.TextFrame.WordWrap = $true

$tr = $discShape.TextFrame.TextRange

# Merge "DISCIPLINA" + ":   " runs into a single "DISCIPLINA:   " run
# (keeps the first run's rPr/formatting).
$run1 = $tr.Characters(1, 10)
$run1.Text = "DISCIPLINA:   "

# Remove the now-redundant ":   " run that followed.
$run2 = $tr.Characters(15, 4)
$run2.Text = ""

# Extend "COMPLIANCE & QUALITY ASSURANCE" with the new suffix text.
$run3 = $tr.Characters(15, 30)
$run3.Text = "COMPLIANCE & QUALITY ASSURANCE  AULA GIT"
